# Fix issues with LocaliserTool
# Updates the SnippetID values (column H) on the "Voice Lines - main" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Voice Lines - main")

# Map of row -> new SnippetID value (column H)
$snippetIds = @{
    2  = "w7zp"
    3  = "w7zp"
    4  = "D5hS"
    5  = "WSyU"
    6  = "aaR0"
    7  = "YDaQ"
    8  = "blNR"
    9  = "GX53"
    10 = "nvdw"
    11 = "SlpG"
    12 = "sd2u"
    13 = "LTD9"
    14 = "wEui"
    15 = "GVje"
    16 = "7inI"
    17 = "y2Bv"
    18 = "y2Bv"
    19 = "y2Bv"
    20 = "y2Bv"
    21 = "y2Bv"
    22 = "hUaM"
    23 = "lAkH"
    24 = "PSMq"
    25 = "xSfR"
    26 = "U06q"
    27 = "U06q"
    28 = "EhGW"
    29 = "ylr1"
}

foreach ($row in $snippetIds.Keys) {
    $ws.Cells.Item($row, 8).Value = $snippetIds[$row]
}
